# Weekly update: insert a new price observation as the latest row for
# "Hortaliza, Agrícola del Norte S.A. de Arica - Ajo".
#
# The new record is inserted at row 8 (pushing the existing rows 8-35
# down to 9-36), mirroring how the source data feed prepends the most
# recent week's observation to the top of the historical series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 8..35 down to 9..36, leaving a blank row 8 to populate.
$ws.Rows.Item(8).Insert()

$ws.Range("A8").Value = 1
$ws.Range("B8").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C8").Value = "Arica y Parinacota"
$ws.Range("D8").Value = 44972
$ws.Range("E8").Value = 15
$ws.Range("F8").Value = 100112003
$ws.Range("G8").Value = "Ajo"
$ws.Range("H8").Value = "Chino"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 550
$ws.Range("K8").Value = 15000
$ws.Range("L8").Value = 16000
$ws.Range("M8").Value = 15636
$ws.Range("N8").Value = "`$/caja 10 kilos"
$ws.Range("O8").Value = "China"
$ws.Range("P8").Value = 1564
$ws.Range("Q8").Value = 10
$ws.Range("R8").Value = "Hortaliza"
